$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 becomes the "Enterprises density (per 1000 people)" row
$ws.Range("A10").Value = "Enterprises density (per 1000 people)"
$ws.Range("C10").ClearContents()

$d10Style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "12.4"
$ws.Range("D10").Style = $d10Style

# Row 11 becomes the "Enterprises (absolute #)" row
$ws.Range("A11").Value = "Enterprises (absolute #)"

$c11Style = $ws.Range("C11").Style
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "55966"
$ws.Range("C11").Style = $c11Style

$d11Style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "55966"
$ws.Range("D11").Style = $d11Style
